# Updates the "cryptos" worksheet with refreshed price/volume figures
# (and, for a few rows, swapped coin identity) as captured by the
# Fri Dec 8 05:13:26 UTC 2023 GitHub Actions data refresh.
#
# The "Price" column (D) sometimes holds digit-only strings (e.g.
# "43.453.62", "0.0000100", "1.00") that Excel would otherwise
# auto-convert into numbers; NumberFormat is forced to Text ('@') for
# those cells before assigning the value so the literal string is kept.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.453.62'
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.369.36'
$ws.Range('E3').Value = '  +4.57%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.62'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.656'
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.39'
$ws.Range('E7').Value = '  +13.46%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.469'
$ws.Range('E9').Value = '  +3.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0977'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.11'
$ws.Range('E11').Value = '  -1.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '27.15'
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.726.58'
$ws.Range('E13').Value = '  +4.76%  '
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.89'
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.31'
$ws.Range('E16').Value = '  +2.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.855'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.375.37'
$ws.Range('E18').Value = '  +5.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.440.25'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000100'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.35'
$ws.Range('E21').Value = '  +3.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.73'
$ws.Range('E22').Value = '  +0.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.97'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.86'
$ws.Range('E24').Value = '  +16.06%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  +1.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.28'
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.88'
$ws.Range('E28').Value = '  +2.03%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.01'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '174.49'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('E31').Value = '  +4.20%  '
$ws.Range('E32').Value = '  -7.77%  '
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0691'
$ws.Range('E35').Value = '  +0.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.08'
$ws.Range('E36').Value = '  +1.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.59'
$ws.Range('E37').Value = '  +1.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.44'
$ws.Range('E38').Value = '  +5.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.67'
$ws.Range('E39').Value = '  -0.64%  '
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '18.75'
$ws.Range('E42').Value = '  +7.64%  '
$ws.Range('E43').Value = '  +1.20%  '
$ws.Range('E44').Value = '  +7.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.53'
$ws.Range('E45').Value = '  +2.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '99.66'
$ws.Range('E46').Value = '  +1.04%  '
$ws.Range('E47').Value = '  +1.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0956'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.441.63'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.596.04'
$ws.Range('E50').Value = '  +4.83%  '
$ws.Range('B51').Value = 'TerraClassic'
$ws.Range('C51').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000204'
$ws.Range('E51').Value = '  -7.81%  '
